$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the comment on C2 (also drops comments1.xml / vmlDrawing / legacyDrawing ref) ---
$ws.Range("C2").Comment.Delete()

# --- Shrink the "run_parameters7" table: drop the "Sample Size" / "Seed Value" columns ---
# Clear the old Sample Size / Seed Value data+formatting (B2:C3) and the soon-to-be-vacated D2:E3,
# then re-populate B2:C3 with what used to live in D2:E3 (Output Path / Version / output / vTest).
$ws.Range("B2:E3").Clear()

$ws.Range("B2").Value = "Output Path"
$ws.Range("C2").Value = "Version"
$ws.Range("B3").Value = "'output"
$ws.Range("C3").Value = "vTest"

$lo1 = $ws.ListObjects.Item(1)
$lo1.Resize($ws.Range("B2:C3"))

# --- New data for the "inputs_from_files" table (rows 7-9) ---
$ws.Range("B7").Value = "File Path"
$ws.Range("C7").Value = "Table Name"
$ws.Range("D7").Value = "Query Only"
$ws.Range("B8").Value = "\adapter\tests\test.db"
$ws.Range("C8").Value = "table1, table2, table3"
$ws.Range("D8").Value = "N, Y, N"
$ws.Range("B9").Value = "\adapter\tests\test.xlsx`t`t"

$lo2 = $ws.ListObjects.Add(1, $ws.Range("B7:D9"), 0, 1)
$ws.ListObjects.Item(1).Name = "inputs_from_files"
$ws.ListObjects.Item(1).TableStyle = "TableStyleMedium6"

# --- Column widths for the new table's columns ---
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.999999999999998
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666

# --- Selection matches the author's final cursor position ---
$null = $ws.Range("H22").Select()
